$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed values for column C (rows 2-9), produced by re-running the
# loop with the corrected logic (see commit message).
$newValues = @{
    2 = -4.998762209006956
    3 = -1.15501353610476
    4 = -0.06001339225539037
    5 = -0.4012163912933333
    6 = 0.01256418049674497
    7 = 0.104613889132365
    8 = 0.1302835008737994
    9 = 0.02670015811043037
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
